# Hungary NB I base update (07-03-2024 23:43):
#   - rows 135-137 (already-played fixtures) get their settlement columns
#     H (FTHG) / I (FTAG) / J (FTR) filled in, plus refreshed closing/live odds
#   - six newly added fixtures become rows 138-143
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Seed rows 138-143 by cloning row 135-137s current layout (A:AC, no result yet) ----
# Do this BEFORE editing rows 135-137 below, since those three rows are about to gain
# H/I/J/AB/AC values and we want the still-unplayed new fixtures to keep the original,
# narrower (A-G,K-AA) shape/format.
$ws.Range("A135:AC135").Copy($ws.Range("A138:AC138"))
$ws.Range("A136:AC136").Copy($ws.Range("A139:AC139"))
$ws.Range("A137:AC137").Copy($ws.Range("A140:AC140"))
$ws.Range("A135:AC135").Copy($ws.Range("A141:AC141"))
$ws.Range("A136:AC136").Copy($ws.Range("A142:AC142"))
$ws.Range("A137:AC137").Copy($ws.Range("A143:AC143"))

# ---- New fixture rows 138-143 ----
# Row 138
$ws.Cells.Item(138,1).Value = 136
$ws.Cells.Item(138,2).Value = 6818337
$ws.Cells.Item(138,3).Value = "Hungary NB I"
$ws.Cells.Item(138,4).Value = "Hungary NB I"
$ws.Cells.Item(138,5).Value = 45359.66666666666
$ws.Cells.Item(138,6).Value = "Puskas Academy"
$ws.Cells.Item(138,7).Value = "Ujpest"
$ws.Cells.Item(138,11).Value = 1.727
$ws.Cells.Item(138,12).Value = 3.6
$ws.Cells.Item(138,13).Value = 4.5
$ws.Cells.Item(138,14).Value = 1.571
$ws.Cells.Item(138,15).Value = 4
$ws.Cells.Item(138,16).Value = 5.5
$ws.Cells.Item(138,17).Value = -1
$ws.Cells.Item(138,18).Value = 2.025
$ws.Cells.Item(138,19).Value = 1.825
$ws.Cells.Item(138,20).Value = 2.5
$ws.Cells.Item(138,21).Value = 1.825
$ws.Cells.Item(138,22).Value = 2.025
$ws.Cells.Item(138,23).Value = 0
$ws.Cells.Item(138,24).Value = 0
$ws.Cells.Item(138,25).Value = 0
$ws.Cells.Item(138,26).Value = 0
$ws.Cells.Item(138,27).Value = 0

# Row 139
$ws.Cells.Item(139,1).Value = 137
$ws.Cells.Item(139,2).Value = 6818336
$ws.Cells.Item(139,3).Value = "Hungary NB I"
$ws.Cells.Item(139,4).Value = "Hungary NB I"
$ws.Cells.Item(139,5).Value = 45360.36458333334
$ws.Cells.Item(139,6).Value = "Mezokovesd Zsory"
$ws.Cells.Item(139,7).Value = "Diosgyori VTK"
$ws.Cells.Item(139,11).Value = 2.9
$ws.Cells.Item(139,12).Value = 3.5
$ws.Cells.Item(139,13).Value = 2.25
$ws.Cells.Item(139,14).Value = 3
$ws.Cells.Item(139,15).Value = 3.5
$ws.Cells.Item(139,16).Value = 2.2
$ws.Cells.Item(139,17).Value = 0.25
$ws.Cells.Item(139,18).Value = 1.85
$ws.Cells.Item(139,19).Value = 2
$ws.Cells.Item(139,20).Value = 2.5
$ws.Cells.Item(139,21).Value = 1.85
$ws.Cells.Item(139,22).Value = 2
$ws.Cells.Item(139,23).Value = 0
$ws.Cells.Item(139,24).Value = 0
$ws.Cells.Item(139,25).Value = 0
$ws.Cells.Item(139,26).Value = 0
$ws.Cells.Item(139,27).Value = 0

# Row 140
$ws.Cells.Item(140,1).Value = 138
$ws.Cells.Item(140,2).Value = 6818339
$ws.Cells.Item(140,3).Value = "Hungary NB I"
$ws.Cells.Item(140,4).Value = "Hungary NB I"
$ws.Cells.Item(140,5).Value = 45360.46875
$ws.Cells.Item(140,6).Value = "Kecskemeti TE"
$ws.Cells.Item(140,7).Value = "Kisvarda FC"
$ws.Cells.Item(140,11).Value = 2.15
$ws.Cells.Item(140,12).Value = 3.2
$ws.Cells.Item(140,13).Value = 3.4
$ws.Cells.Item(140,14).Value = 2.1
$ws.Cells.Item(140,15).Value = 3.2
$ws.Cells.Item(140,16).Value = 3.6
$ws.Cells.Item(140,17).Value = -0.25
$ws.Cells.Item(140,18).Value = 1.875
$ws.Cells.Item(140,19).Value = 1.975
$ws.Cells.Item(140,20).Value = 2.25
$ws.Cells.Item(140,21).Value = 1.825
$ws.Cells.Item(140,22).Value = 2.025
$ws.Cells.Item(140,23).Value = 0
$ws.Cells.Item(140,24).Value = 0
$ws.Cells.Item(140,25).Value = 0
$ws.Cells.Item(140,26).Value = 0
$ws.Cells.Item(140,27).Value = 0

# Row 141
$ws.Cells.Item(141,1).Value = 139
$ws.Cells.Item(141,2).Value = 6818335
$ws.Cells.Item(141,3).Value = "Hungary NB I"
$ws.Cells.Item(141,4).Value = "Hungary NB I"
$ws.Cells.Item(141,5).Value = 45361.375
$ws.Cells.Item(141,6).Value = "Paksi"
$ws.Cells.Item(141,7).Value = "Zalaegerszegi TE"
$ws.Cells.Item(141,11).Value = 2.1
$ws.Cells.Item(141,12).Value = 3.4
$ws.Cells.Item(141,13).Value = 3.4
$ws.Cells.Item(141,14).Value = 2
$ws.Cells.Item(141,15).Value = 3.4
$ws.Cells.Item(141,16).Value = 3.6
$ws.Cells.Item(141,17).Value = -0.5
$ws.Cells.Item(141,18).Value = 2.025
$ws.Cells.Item(141,19).Value = 1.825
$ws.Cells.Item(141,20).Value = 2.5
$ws.Cells.Item(141,21).Value = 1.825
$ws.Cells.Item(141,22).Value = 2.025
$ws.Cells.Item(141,23).Value = 0
$ws.Cells.Item(141,24).Value = 0
$ws.Cells.Item(141,25).Value = 0
$ws.Cells.Item(141,26).Value = 0
$ws.Cells.Item(141,27).Value = 0

# Row 142
$ws.Cells.Item(142,1).Value = 140
$ws.Cells.Item(142,2).Value = 6818338
$ws.Cells.Item(142,3).Value = "Hungary NB I"
$ws.Cells.Item(142,4).Value = "Hungary NB I"
$ws.Cells.Item(142,5).Value = 45361.47916666666
$ws.Cells.Item(142,6).Value = "MOL Fehervar FC"
$ws.Cells.Item(142,7).Value = "Ferencvarosi TC"
$ws.Cells.Item(142,11).Value = 5.25
$ws.Cells.Item(142,12).Value = 4
$ws.Cells.Item(142,13).Value = 1.571
$ws.Cells.Item(142,14).Value = 5.5
$ws.Cells.Item(142,15).Value = 4.2
$ws.Cells.Item(142,16).Value = 1.533
$ws.Cells.Item(142,17).Value = 1
$ws.Cells.Item(142,18).Value = 1.975
$ws.Cells.Item(142,19).Value = 1.875
$ws.Cells.Item(142,20).Value = 3
$ws.Cells.Item(142,21).Value = 2.05
$ws.Cells.Item(142,22).Value = 1.8
$ws.Cells.Item(142,23).Value = 0
$ws.Cells.Item(142,24).Value = 0
$ws.Cells.Item(142,25).Value = 0
$ws.Cells.Item(142,26).Value = 0
$ws.Cells.Item(142,27).Value = 0

# Row 143
$ws.Cells.Item(143,1).Value = 141
$ws.Cells.Item(143,2).Value = 6818334
$ws.Cells.Item(143,3).Value = "Hungary NB I"
$ws.Cells.Item(143,4).Value = "Hungary NB I"
$ws.Cells.Item(143,5).Value = 45361.66666666666
$ws.Cells.Item(143,6).Value = "Debreceni VSC"
$ws.Cells.Item(143,7).Value = "MTK Budapest"
$ws.Cells.Item(143,11).Value = 2.2
$ws.Cells.Item(143,12).Value = 3.4
$ws.Cells.Item(143,13).Value = 3.1
$ws.Cells.Item(143,14).Value = 2.2
$ws.Cells.Item(143,15).Value = 3.4
$ws.Cells.Item(143,16).Value = 3.2
$ws.Cells.Item(143,17).Value = -0.25
$ws.Cells.Item(143,18).Value = 1.925
$ws.Cells.Item(143,19).Value = 1.925
$ws.Cells.Item(143,20).Value = 2.5
$ws.Cells.Item(143,21).Value = 1.825
$ws.Cells.Item(143,22).Value = 2.025
$ws.Cells.Item(143,23).Value = 0
$ws.Cells.Item(143,24).Value = 0
$ws.Cells.Item(143,25).Value = 0
$ws.Cells.Item(143,26).Value = 0
$ws.Cells.Item(143,27).Value = 0

# ---- Existing rows 135-137: fill in result columns + refresh odds ----
# Row 135
$ws.Cells.Item(135,8).Value = 0
$ws.Cells.Item(135,9).Value = 2
$ws.Cells.Item(135,10).Value = "A"
$ws.Cells.Item(135,14).Value = 3.2
$ws.Cells.Item(135,16).Value = 2
$ws.Cells.Item(135,18).Value = 1.8
$ws.Cells.Item(135,19).Value = 2.05
$ws.Cells.Item(135,21).Value = 1.925
$ws.Cells.Item(135,22).Value = 1.925
$ws.Cells.Item(135,23).Value = -1
$ws.Cells.Item(135,24).Value = -1
$ws.Cells.Item(135,25).Value = 1
$ws.Cells.Item(135,26).Value = -1
$ws.Cells.Item(135,27).Value = 1.05
$ws.Cells.Item(135,28).Value = -1
$ws.Cells.Item(135,29).Value = 0.925

# Row 136
$ws.Cells.Item(136,8).Value = 2
$ws.Cells.Item(136,9).Value = 0
$ws.Cells.Item(136,10).Value = "H"
$ws.Cells.Item(136,14).Value = 1.222
$ws.Cells.Item(136,15).Value = 6.5
$ws.Cells.Item(136,16).Value = 8
$ws.Cells.Item(136,17).Value = -1.75
$ws.Cells.Item(136,21).Value = 1.95
$ws.Cells.Item(136,22).Value = 1.9
$ws.Cells.Item(136,23).Value = 0.222
$ws.Cells.Item(136,24).Value = -1
$ws.Cells.Item(136,25).Value = -1
$ws.Cells.Item(136,26).Value = 0.475
$ws.Cells.Item(136,27).Value = -0.5
$ws.Cells.Item(136,28).Value = -1
$ws.Cells.Item(136,29).Value = 0.8999999999999999

# Row 137
$ws.Cells.Item(137,8).Value = 2
$ws.Cells.Item(137,9).Value = 0
$ws.Cells.Item(137,10).Value = "H"
$ws.Cells.Item(137,14).Value = 2.1
$ws.Cells.Item(137,15).Value = 3.4
$ws.Cells.Item(137,16).Value = 3.2
$ws.Cells.Item(137,17).Value = -0.25
$ws.Cells.Item(137,18).Value = 1.85
$ws.Cells.Item(137,19).Value = 2
$ws.Cells.Item(137,21).Value = 2.025
$ws.Cells.Item(137,22).Value = 1.825
$ws.Cells.Item(137,23).Value = 1.1
$ws.Cells.Item(137,24).Value = -1
$ws.Cells.Item(137,25).Value = -1
$ws.Cells.Item(137,26).Value = 0.8500000000000001
$ws.Cells.Item(137,27).Value = -1
$ws.Cells.Item(137,28).Value = -1
$ws.Cells.Item(137,29).Value = 0.825
